# Texas COVID-19 New Confirmed Cases by County
# Add the three newest daily "New Cases" columns (12-01, 12-02, 12-03),
# their values for Collin county, and tidy up the sheet view, matching
# the "Fully implemented all data for the today query" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new date-header columns (row 1) ---
$ws.Range("JK1").Value = "New Cases 12-01"
$ws.Range("JL1").Value = "New Cases 12-02"
$ws.Range("JM1").Value = "New Cases 12-03"

# --- new daily counts for Collin county (row 2) ---
$ws.Range("JK2").Value = 390
$ws.Range("JL2").Value = 561
$ws.Range("JM2").Value = 295

# match the bold/centered/bordered header style used by the rest of row 1
$ws.Range("JJ1").Copy()
$ws.Range("JK1:JM1").PasteSpecial(-4122)

# normalize column widths across the now-wider used range
$ws.Range("A1:JM1").EntireColumn.ColumnWidth = 11.1666666666667

# drop the frozen header pane/old multi-pane selection, leaving a simple
# single selection like the saved workbook
$excel.ActiveWindow.FreezePanes = $false
[void]$ws.Range("E11").Select()
